$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This string is used by the per-language "Status" columns on every sheet
# (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4). Update every occurrence so the
# shared string collapses cleanly to the new text everywhere it appears.

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Range("E4").Value = "In Translation"
$ws1.Range("F4").Value = "In Translation"

$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item(3)   # de-de
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("C4").Value = "In Translation"

# --- Column width change: narrow the Status columns ---
# Old stored width 17.2159881591797 chars -> new stored width 13.4101845877511 chars.
# Range.ColumnWidth is quantized by the host to the nearest 1/6 of a character,
# so we pick the input (12.5) whose rounded result (13.333333...) lands closest
# to the target width.
$ws1.Range("E1").ColumnWidth = 12.5
$ws1.Range("F1").ColumnWidth = 12.5
$ws2.Range("C1").ColumnWidth = 12.5
$ws3.Range("C1").ColumnWidth = 12.5
